$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.153.96'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +3.82%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.105.73'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.71%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '561.81'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.38'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +6.71%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.102.04'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.69%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.501'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.95%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.16'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +15.54%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.154'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.57%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.469'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +3.77%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000234'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +4.52%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.56'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.86%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.603.08'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.62%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.074.70'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.69%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.102.38'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.56%  '

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.48%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.89'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.35%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '483.46'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.53%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.87'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.47%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.682'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.11%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.56'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +6.29%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.53'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +10.69%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.40'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.15%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.03%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.80'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.38%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.24'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +4.24%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.07'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +5.79%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.12%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.20'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.78%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.71%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.50'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +5.23%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.70'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.91%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.26'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +5.39%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.27'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.30%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '473.13'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.62%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0412'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +5.49%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0834'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +3.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.95'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +19.54%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.019.32'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.32'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.56%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.117'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.76%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '28.36'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +6.13%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.259'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +4.60%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.14'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +7.24%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.54%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₃0525'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +5.47%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '117.16'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.35%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.09'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.55%  '
